# Refresh cryptocurrency price / 1h-volume change figures (and fix the
# FraxShare / TrustWalletToken row ordering) to match the latest scrape.
# Values that look like plain numbers are entered with a leading single
# quote so Excel stores them as text (matching the sheet's inlineStr cells)
# instead of silently converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.759.69'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '2.104.50'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Value = "'" + '0.617'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = "'" + '62.28'
$ws.Range("E7").Value = '  +1.18%  '
$ws.Range("D8").Value = "'" + '1.00'
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("D10").Value = "'" + '0.0843'
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("D12").Value = "'" + '15.78'
$ws.Range("E12").Value = '  +6.41%  '
$ws.Range("D13").Value = '2.416.87'
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = "'" + '22.16'
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '2.104.19'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '38.802.34'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").Value = "'" + '71.98'
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("D20").Value = "'" + '6.11'
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").Value = "'" + '228.15'
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = "'" + '2.34'
$ws.Range("E24").Value = '  -4.23%  '
$ws.Range("D25").Value = "'" + '2.31'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").Value = "'" + '171.73'
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("D27").Value = "'" + '9.56'
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("E28").Value = '  +5.15%  '
$ws.Range("E29").Value = '  +4.09%  '
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("E31").Value = '  +7.69%  '
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("E33").Value = '  +1.54%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  +7.56%  '
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("E37").Value = '  +1.18%  '
$ws.Range("D38").Value = "'" + '3.59'
$ws.Range("E38").Value = '  +1.73%  '
$ws.Range("D39").Value = "'" + '1.00'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").Value = "'" + '18.11'
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("D41").Value = "'" + '102.82'
$ws.Range("E41").Value = '  +2.50%  '
$ws.Range("E42").Value = '  +3.40%  '
$ws.Range("D43").Value = '1.537.72'
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = "'" + '7.88'
$ws.Range("E44").Value = '  +4.28%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = "'" + '1.16'
$ws.Range("E45").Value = '  +4.06%  '
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("D51").Value = '2.303.23'
$ws.Range("E51").Value = '  +0.19%  '
